# Actualizacion automatica 2025-10-15 10:30:10
# A new client ("GAVILANES VELEZ MARIA VALERIA") was added to the sales
# roster. This inserts one row before the existing alphabetical block
# (between "GARCIA BRAVO JOSE LUIS" and "GRANJA VANEGAS MARCELA") on the
# two detail sheets, shifting all subsequent client rows down by one, and
# refreshes the derived "de 49" -> "de 50" counters plus the roll-up
# figures on the summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"  (A1:R51 -> A1:R52)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Push rows 24..51 down to 25..52, leaving a blank row 24 behind.
$ws1.Rows.Item(24).Insert()

# Populate the newly-opened row 24 with the new client (all zero sales).
$ws1.Range("A24").Value = "OFICINA-CATAECSA"
$ws1.Range("B24").Value = "GAVILANES VELEZ MARIA VALERIA"
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(24, $c).Value = 0
}

# The final counter row (old row 51, now row 52) reports "x de 49" per
# column; with one more client row it becomes "x de 50".
for ($c = 3; $c -le 18; $c++) {
    $cell = $ws1.Cells.Item(52, $c)
    $cell.Value = $cell.Value2.Replace("de 49", "de 50")
}

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"  (A1:G51 -> A1:G52)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(24).Insert()

$ws2.Range("A24").Value = "OFICINA-CATAECSA"
$ws2.Range("B24").Value = "GAVILANES VELEZ MARIA VALERIA"
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(24, $c).Value = 0
}
# Row 52 keeps the same totals (2652.23 / 1566.5 / 13412.01 / 574.9 / 0) -
# the new client contributes nothing, so no further edit is needed there.

# ---------------------------------------------------------------------
# Sheet 3: "CUMPLIMIENTO MENSUAL"  (monthly roll-up, values only)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D3").Value = 5264.61
$ws3.Range("E3").Value = 14735.39
$ws3.Range("F3").Value = 0.2632305

$ws3.Range("D4").Value = 8611.43
$ws3.Range("E4").Value = 11388.57
$ws3.Range("F4").Value = 0.4305715
